$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Story line Rockville block (originally lettered 5a-5e / 6f)
# ---------------------------------------------------------------------------

# Paragraph currently reading "5a.  if relic equals yes..."
$p = $d.Paragraphs(21)

# 1) Insert a brand-new bullet BEFORE it: "5a. set user input to all lower case"
$p.Range.InsertParagraphBefore()
$d.Paragraphs(21).Range.Text = "5a. set user input to all lower case"

# 2) The old "5a.  if relic equals yes..." paragraph (now index 22) becomes "5b. ..."
$target = $d.Paragraphs(22).Range
$target.Find.Execute("a.  if relic equals yes", $true, $false, $false, $false, $false, $true, 1, $false, "b.  if relic equals yes", 1)

# 3) "5b. if dollar amount is < 50 ..." -> "5c. if dollar amount is < 50 ..."
$target = $d.Paragraphs(23).Range
$target.Find.Execute("b", $true, $false, $false, $false, $false, $true, 1, $false, "c", 1)

# 4) "5c.Then Output ..." -> "5d.Then Output ..." (letter sits between proofErr gramStart/gramEnd)
$target = $d.Paragraphs(24).Range
$target.Find.Execute("c", $true, $false, $false, $false, $false, $true, 1, $false, "d", 1)

# 5) "5d. if dollar amount >= 50 ..." -> "5e. if dollar amount >= 50 ..."
$target = $d.Paragraphs(25).Range
$target.Find.Execute("d", $true, $false, $false, $false, $false, $true, 1, $false, "e", 1)

# 6) "5e.Then output 'GAME OVER: You went home broke'" -> "5f.Then output ..." (proofErr bounded)
$target = $d.Paragraphs(26).Range
$target.Find.Execute("e", $true, $false, $false, $false, $false, $true, 1, $false, "f", 1)

# 7) "6f. then output 'GAME OVER: you went home empty handed'" -> "6a. then output ..."
$target = $d.Paragraphs(28).Range
$target.Find.Execute("f", $true, $false, $false, $false, $false, $true, 1, $false, "a", 1)

# ---------------------------------------------------------------------------
# Story line Forest block (originally lettered 7a-7h)
# ---------------------------------------------------------------------------

# Paragraph currently reading "7f. if reward equal gold or rubies..."
$p2 = $d.Paragraphs(37)

# 8) Insert a brand-new bullet BEFORE it: "7f. set treasure input to all lowercase"
$p2.Range.InsertParagraphBefore()
$d.Paragraphs(37).Range.Text = "7f. set treasure input to all lowercase"

# 9) The old "7f. if reward equal gold or rubies..." paragraph (now index 38) becomes "7g. ..."
$target = $d.Paragraphs(38).Range
$target.Find.Execute("f", $true, $false, $false, $false, $false, $true, 1, $false, "g", 1)

# 10) "7f. if reward equal a ride home..." -> "7I. if reward equal a ride home..."
$target = $d.Paragraphs(40).Range
$target.Find.Execute("f", $true, $false, $false, $false, $false, $true, 1, $false, "I", 1)

# 11) "7g. then output 'GAME OVER: you got home ...'" -> "7J. then output ..."
$target = $d.Paragraphs(41).Range
$target.Find.Execute("g", $true, $false, $false, $false, $false, $true, 1, $false, "J", 1)
